$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected naive-forecaster error values (rows 24-52, columns B:K)
# Maps cell address -> corrected value, applied via Range(...).Value
$updates = @{
    "K24" = -8.183026219731815
    "J25" = -8.028721684813732
    "K25" = 0.1109553426839227
    "I26" = -8.680305474188014
    "J26" = -0.54062844669036
    "K26" = 1.60694255667963
    "H27" = -8.462984671940774
    "I27" = -0.3233076444431192
    "J27" = 1.824263358926871
    "K27" = -1.615771047272361
    "G28" = -8.080305474188014
    "H28" = 0.05937155330964
    "I28" = 2.20694255667963
    "J28" = -1.233091849519602
    "K28" = -1.294275412034003
    "F29" = -8.554988394376949
    "G29" = -0.4153113668792941
    "H29" = 1.732259636490696
    "I29" = -1.707774769708536
    "J29" = -1.768958332222937
    "K29" = 0.2725059170203678
    "E30" = -8.396144372968052
    "F30" = -0.2564673454703979
    "G30" = 1.891103657899592
    "H30" = -1.54893074829964
    "I30" = -1.610114310814041
    "J30" = 0.431349938429264
    "K30" = -0.1512409947814746
    "D31" = -8.452313642352863
    "E31" = -0.3126366148552081
    "F31" = 1.834934388514782
    "G31" = -1.60510001768445
    "H31" = -1.666283580198851
    "I31" = 0.3751806690444539
    "J31" = -0.2074102641662847
    "K31" = -0.1117330153097526
    "C32" = -10.68030547418801
    "D32" = -2.54062844669036
    "E32" = -0.3930574433203702
    "F32" = -3.833091849519602
    "G32" = -3.894275412034003
    "H32" = -1.852811162790698
    "I32" = -2.435402096001437
    "J32" = -2.339724847144905
    "K32" = -2.60014685520589
    "B33" = -18.17126098581633
    "C33" = -10.03158395831868
    "D33" = -7.884012954948687
    "E33" = -11.32404736114792
    "F33" = -11.38523092366232
    "G33" = -9.343766674419015
    "H33" = -9.926357607629754
    "I33" = -9.830680358773222
    "J33" = -10.09110236683421
    "K33" = -9.255644381000664
    "B34" = 8.139677027497655
    "C34" = 10.28724803086764
    "D34" = 6.847213624668413
    "E34" = 6.786030062154012
    "F34" = 8.827494311397317
    "G34" = 8.244903378186578
    "H34" = 8.34058062704311
    "I34" = 8.080158618982125
    "J34" = 8.915616604815668
    "K34" = 8.81559383571846
    "B35" = 2.14757100336999
    "C35" = -1.292463402829242
    "D35" = -1.353646965343643
    "E35" = 0.6878172838996619
    "F35" = 0.1052263506889233
    "G35" = 0.2009035995454554
    "H35" = -0.05951840851552959
    "I35" = 0.7759395773180131
    "J35" = 0.6759168082208049
    "K35" = 0.3220587481868762
    "B36" = -3.440034406199232
    "C36" = -3.501217968713633
    "D36" = -1.459753719470328
    "E36" = -2.042344652681066
    "F36" = -1.946667403824535
    "G36" = -2.20708941188552
    "H36" = -1.371631426051977
    "I36" = -1.471654195149185
    "J36" = -1.825512255183114
    "K36" = -1.677319185747749
    "B37" = -0.06118356251440082
    "C37" = 1.980280686728904
    "D37" = 1.397689753518165
    "E37" = 1.493367002374697
    "F37" = 1.232944994313713
    "G37" = 2.068402980147255
    "H37" = 1.968380211050047
    "I37" = 1.614522151016118
    "J37" = 1.762715220451483
    "K37" = 1.920148881959748
    "B38" = 2.041464249243305
    "C38" = 1.458873316032566
    "D38" = 1.554550564889098
    "E38" = 1.294128556828113
    "F38" = 2.129586542661656
    "G38" = 2.029563773564448
    "H38" = 1.675705713530519
    "I38" = 1.823898782965884
    "J38" = 1.981332444474148
    "K38" = 1.480676358197826
    "B39" = -0.5825909332107386
    "C39" = -0.4869136843542065
    "D39" = -0.7473356924151915
    "E39" = 0.08812229341835121
    "F39" = -0.01190047567885699
    "G39" = -0.3657585357127857
    "H39" = -0.2175654662774206
    "I39" = -0.06013180476915631
    "J39" = -0.5607878910454785
    "K39" = -0.2804316037748719
    "B40" = 0.09567724885653209
    "C40" = -0.1647447592044529
    "D40" = 0.6707132266290898
    "E40" = 0.5706904575318816
    "F40" = 0.2168323974979529
    "G40" = 0.365025466933318
    "H40" = 0.5224591284415823
    "I40" = 0.02180304216525999
    "J40" = 0.3021593294358667
    "K40" = 0.1301932178764815
    "B41" = -0.260422008060985
    "C41" = 0.5750359777725577
    "D41" = 0.4750132086753495
    "E41" = 0.1211551486414208
    "F41" = 0.2693482180767859
    "G41" = 0.4267818795850502
    "H41" = -0.07387420669127209
    "I41" = 0.2064820805793346
    "J41" = 0.0345159690199494
    "K41" = 0.3402630873498395
    "B42" = 0.8354579858335427
    "C42" = 0.7354352167363345
    "D42" = 0.3815771567024058
    "E42" = 0.5297702261377709
    "F42" = 0.6872038876460351
    "G42" = 0.1865478013697129
    "H42" = 0.4669040886403196
    "I42" = 0.2949379770809344
    "J42" = 0.6006850954108245
    "K42" = -0.01150507022163078
    "B43" = -0.1000227690972082
    "C43" = -0.4538808291311369
    "D43" = -0.3056877596957718
    "E43" = -0.1482540981875075
    "F43" = -0.6489101844638298
    "G43" = -0.3685538971932231
    "H43" = -0.5405200087526083
    "I43" = -0.2347728904227182
    "J43" = -0.8469630560551735
    "K43" = -0.1586268237156929
    "B44" = -0.3538580600339287
    "C44" = -0.2056649905985636
    "D44" = -0.04823132909029931
    "E44" = -0.5488874153666216
    "F44" = -0.2685311280960149
    "G44" = -0.4404972396554001
    "H44" = -0.13475012132551
    "I44" = -0.7469402869579653
    "J44" = -0.05860405461848467
    "B45" = 0.1481930694353651
    "C45" = 0.3056267309436294
    "D45" = -0.1950293553326929
    "E45" = 0.08532693193791374
    "F45" = -0.08663917962147143
    "G45" = 0.2191079387084187
    "H45" = -0.3930822269240366
    "I45" = 0.295254005415444
    "B46" = 0.1574336615082643
    "C46" = -0.343222424768058
    "D46" = -0.06286613749745135
    "E46" = -0.2348322490568365
    "F46" = 0.0709148692730536
    "G46" = -0.5412752963594016
    "H46" = 0.1470609359800789
    "B47" = -0.5006560862763223
    "C47" = -0.2202997990057156
    "D47" = -0.3922659105651008
    "E47" = -0.08651879223521067
    "F47" = -0.698708957867666
    "G47" = -0.01037272552818536
    "B48" = 0.2803562872706067
    "C48" = 0.1083901757112215
    "D48" = 0.4141372940411116
    "E48" = -0.1980528715913437
    "F48" = 0.4902833607481369
    "B49" = -0.1719661115593852
    "C49" = 0.1337810067705049
    "D49" = -0.4784091588619503
    "E49" = 0.2099270734775303
    "B50" = 0.3057471183298901
    "C50" = -0.3064430473025652
    "D50" = 0.3818931850369154
    "B51" = -0.6121901656324553
    "C51" = 0.07614606670702531
    "B52" = 0.6883362323394806
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
